$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PayNowCorpNoCF")
$ws.Range("B2").Value = "Fri Nov 21 00:37:24 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCorpNoCFOnly")
$ws.Range("B2").Value = "Fri Nov 21 00:40:53 IST 2025"

$ws = $wb.Worksheets.Item("ACMismatchCorp")
$ws.Range("B2").Value = "Thu Nov 20 22:50:49 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCorpDCF")
$ws.Range("B2").Value = "Fri Nov 21 00:30:56 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCorpSCF")
$ws.Range("B2").Value = "Fri Nov 21 00:43:26 IST 2025"

$ws = $wb.Worksheets.Item("MaxAmountErrorCorp")
$ws.Range("B2").Value = "Thu Nov 20 22:53:33 IST 2025"

$ws = $wb.Worksheets.Item("MinAmountErrorCorp")
$ws.Range("B2").Value = "Thu Nov 20 22:55:37 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Thu Nov 20 23:41:24 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCorpNoCFReqFields")
$ws.Range("B2").Value = "Fri Nov 21 00:42:55 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPCNoCF")
$ws.Range("B2").Value = "Fri Nov 21 00:47:57 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPCSCF")
$ws.Range("B2").Value = "Fri Nov 21 01:01:12 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPCDCF")
$ws.Range("B2").Value = "Fri Nov 21 00:45:53 IST 2025"

$ws = $wb.Worksheets.Item("ACMismatchPC")
$ws.Range("B2").Value = "Thu Nov 20 22:51:41 IST 2025"

$ws = $wb.Worksheets.Item("DualCFCeilingCorp")
$ws.Range("B2").Value = "Thu Nov 20 23:11:40 IST 2025"

$ws = $wb.Worksheets.Item("DualCFCeilingPC")
$ws.Range("B2").Value = "Thu Nov 20 23:14:27 IST 2025"

$ws = $wb.Worksheets.Item("DualCFCeilingPS")
$ws.Range("B2").Value = "Thu Nov 20 23:15:03 IST 2025"

$ws = $wb.Worksheets.Item("DualCFFlatCorp")
$ws.Range("B2").Value = "Thu Nov 20 23:19:52 IST 2025"

$ws = $wb.Worksheets.Item("DualCFFlatPC")
$ws.Range("B2").Value = "Thu Nov 20 23:20:33 IST 2025"

$ws = $wb.Worksheets.Item("DualCFFlatPS")
$ws.Range("B2").Value = "Thu Nov 20 23:21:15 IST 2025"

$ws = $wb.Worksheets.Item("DualCFPercentageCorp")
$ws.Range("B2").Value = "Thu Nov 20 23:22:35 IST 2025"

$ws = $wb.Worksheets.Item("DualCFPercentagePC")
$ws.Range("B2").Value = "Thu Nov 20 23:25:04 IST 2025"

$ws = $wb.Worksheets.Item("DualCFPercentagePS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Thu Nov 20 23:32:59 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFCeilingCorp")
$ws.Range("B2").Value = "Thu Nov 20 23:53:27 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFCeilingPS")
$ws.Range("B2").Value = "Thu Nov 20 23:54:39 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFCeilingPC")
$ws.Range("B2").Value = "Thu Nov 20 23:54:06 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFFlatCorp")
$ws.Range("B2").Value = "Thu Nov 20 23:57:43 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFFlatPC")
$ws.Range("B2").Value = "Thu Nov 20 23:58:25 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFFlatPS")
$ws.Range("B2").Value = "Fri Nov 21 00:00:58 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFPercentageCorp")
$ws.Range("B2").Value = "Fri Nov 21 00:07:40 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFPercentagePC")
$ws.Range("B2").Value = "Fri Nov 21 00:14:36 IST 2025"

$ws = $wb.Worksheets.Item("SingleCFPercentagePS")
$ws.Range("B2").Value = "Fri Nov 21 00:20:13 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPCNoCFReqFields")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Nov 21 01:00:43 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPCNoCFOnly")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Fri Nov 21 00:54:26 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountPC")
$ws.Range("B2").Value = "Thu Nov 20 23:46:27 IST 2025"

$ws = $wb.Worksheets.Item("MinAmountErrorPC")
$ws.Range("B2").Value = "Thu Nov 20 22:56:02 IST 2025"

$ws = $wb.Worksheets.Item("MaxAmountErrorPC")
$ws.Range("B2").Value = "Thu Nov 20 22:54:12 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPSNoCF")
$ws.Range("B2").Value = "Fri Nov 21 01:05:13 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPSSCF")
$ws.Range("B2").Value = "Fri Nov 21 01:13:49 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPSDCF")
$ws.Range("B2").Value = "Fri Nov 21 01:04:28 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPSNoCFOnly")
$ws.Range("B2").Value = "Fri Nov 21 01:07:41 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPSNoCFReqFields")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Nov 21 01:13:15 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountPS")
$ws.Range("B2").Value = "Thu Nov 20 23:48:34 IST 2025"

$ws = $wb.Worksheets.Item("MaxAmountErrorPS")
$ws.Range("B2").Value = "Thu Nov 20 22:54:40 IST 2025"

$ws = $wb.Worksheets.Item("MinAmountErrorPS")
$ws.Range("B2").Value = "Thu Nov 20 22:56:26 IST 2025"
